# "added scripts for event annotation"
# Populate a new column G ("Unrelated") on every data row of Sheet1 except
# the handful of rows that already carried a value in G (26, 36, 65, 80),
# then update the sheet selection and the width of column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that already have data in column G and must be left untouched.
$skipRows = @(26, 36, 65, 80)

for ($r = 2; $r -le 92; $r++) {
    if ($skipRows -contains $r) {
        continue
    }
    $ws.Cells.Item($r, 7).Value = "Unrelated"
}

# Narrow column E (was ~36 chars wide) down to ~20 chars wide.
$ws.Columns.Item(5).ColumnWidth = 19.25

# Update the visible selection to match the edited range.
$ws.Range("G81:G92").Select()
